$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2325.2307
$ws.Range("I40").Value = 2432
$ws.Range("J40").Value = 2085
$ws.Range("K40").Value = 2432
$ws.Range("L40").Value = 2085
$ws.Range("M40").Value = -2257
$ws.Range("N40").Value = -2435
$ws.Range("H43").Value = 10589.1875
$ws.Range("I43").Value = 15383.167
$ws.Range("K43").Value = 15383.167
$ws.Range("M43").Value = -15314.167
$ws.Range("H138").Value = 2540.487
$ws.Range("I138").Value = 1800.6428
$ws.Range("J138").Value = 2954.8
$ws.Range("K138").Value = 5401.928400000001
$ws.Range("L138").Value = 8864.400000000001
$ws.Range("M138").Value = -261.9284000000007
$ws.Range("N138").Value = -19144.4

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21261.371
$ws.Range("I32").Value = 22754.18
$ws.Range("J32").Value = 2601.25
$ws.Range("K32").Value = 22754.18
$ws.Range("L32").Value = 2601.25
$ws.Range("M32").Value = -22467.18
$ws.Range("N32").Value = -3175.25
$ws.Range("H61").Value = 6806.125
$ws.Range("I61").Value = 5957.737
$ws.Range("K61").Value = 5957.737
$ws.Range("M61").Value = -5745.737
$ws.Range("H136").Value = 6806.125
$ws.Range("I136").Value = 5957.737
$ws.Range("K136").Value = 17873.211
$ws.Range("M136").Value = -15323.211

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 904
$ws.Range("I37").Value = 904
$ws.Range("K37").Value = 904
$ws.Range("M37").Value = -767
$ws.Range("H105").Value = 55571404
$ws.Range("I105").Value = 55571404
$ws.Range("K105").Value = 55571404
$ws.Range("M105").Value = -55569657
$ws.Range("H107").Value = 1361.8
$ws.Range("I107").Value = 1373.8276
$ws.Range("K107").Value = 1373.8276
$ws.Range("M107").Value = 546.1723999999999
$ws.Range("H134").Value = 20516
$ws.Range("I134").Value = 30637.8
$ws.Range("J134").Value = 12081.167
$ws.Range("K134").Value = 91913.39999999999
$ws.Range("L134").Value = 36243.501
$ws.Range("M134").Value = -89378.39999999999
$ws.Range("N134").Value = -41313.501

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1416.85
$ws.Range("I22").Value = 447.92307
$ws.Range("K22").Value = 447.92307
$ws.Range("M22").Value = -97.92307
$ws.Range("H35").Value = 689.5294
$ws.Range("I35").Value = 794.0769
$ws.Range("J35").Value = 349.75
$ws.Range("K35").Value = 794.0769
$ws.Range("L35").Value = 349.75
$ws.Range("M35").Value = -500.0769
$ws.Range("N35").Value = -937.75
$ws.Range("H58").Value = 5338.2607
$ws.Range("I58").Value = 2804.1428
$ws.Range("J58").Value = 9280.223
$ws.Range("K58").Value = 2804.1428
$ws.Range("L58").Value = 9280.223
$ws.Range("M58").Value = -2601.1428
$ws.Range("N58").Value = -9686.223
$ws.Range("H132").Value = 34708.773
$ws.Range("I132").Value = 3383.0688
$ws.Range("J132").Value = 117294.73
$ws.Range("K132").Value = 10149.2064
$ws.Range("L132").Value = 351884.19
$ws.Range("M132").Value = -7619.206399999999
$ws.Range("N132").Value = -356944.19
$ws.Range("H134").Value = 5546.606
$ws.Range("I134").Value = 5819.6294
$ws.Range("J134").Value = 4318
$ws.Range("K134").Value = 17458.8882
$ws.Range("L134").Value = 12954
$ws.Range("M134").Value = -14923.8882
$ws.Range("N134").Value = -18024
$ws.Range("H136").Value = 5338.2607
$ws.Range("I136").Value = 2804.1428
$ws.Range("J136").Value = 9280.223
$ws.Range("K136").Value = 8412.428400000001
$ws.Range("L136").Value = 27840.669
$ws.Range("M136").Value = -5862.428400000001
$ws.Range("N136").Value = -32940.669

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35231840
$ws.Range("I4").Value = 50319170
$ws.Range("J4").Value = 6428761
$ws.Range("K4").Value = 150957510
$ws.Range("L4").Value = 19286283
$ws.Range("M4").Value = -150957398
$ws.Range("N4").Value = -19286507
$ws.Range("H37").Value = 163774
$ws.Range("J37").Value = 163774
$ws.Range("L37").Value = 491322
$ws.Range("N37").Value = -491546
$ws.Range("H50").Value = 801.7727
$ws.Range("J50").Value = 1162.4166
$ws.Range("L50").Value = 3487.2498
$ws.Range("N50").Value = -4449.2498
$ws.Range("H53").Value = 801.7727
$ws.Range("J53").Value = 1162.4166
$ws.Range("L53").Value = 3487.2498
$ws.Range("N53").Value = -4449.2498
$ws.Range("H107").Value = 655.25
$ws.Range("I107").Value = 302
$ws.Range("J107").Value = 712.2258
$ws.Range("K107").Value = 906
$ws.Range("L107").Value = 2136.6774
$ws.Range("M107").Value = 1014
$ws.Range("N107").Value = -5976.6774

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15970
$ws.Range("H97").Value = 1340.4615
$ws.Range("I97").Value = 1407.6
$ws.Range("J97").Value = 1116.6666
$ws.Range("K97").Value = 1407.6
$ws.Range("L97").Value = 1116.6666
$ws.Range("M97").Value = -911.5999999999999
$ws.Range("N97").Value = -2108.6666
$ws.Range("H132").Value = 11809.4
$ws.Range("I132").Value = 8399.4
$ws.Range("J132").Value = 15219.4
$ws.Range("K132").Value = 25198.2
$ws.Range("L132").Value = 45658.2
$ws.Range("M132").Value = -22668.2
$ws.Range("N132").Value = -50718.2

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 202399.8
$ws.Range("I7").Value = 335099.66
$ws.Range("K7").Value = 335099.66
$ws.Range("M7").Value = -334987.66
$ws.Range("H126").Value = 202399.8
$ws.Range("I126").Value = 335099.66
$ws.Range("K126").Value = 1005298.98
$ws.Range("M126").Value = -1002828.98
$ws.Range("H132").Value = 8431.861999999999
$ws.Range("I132").Value = 7596.7617
$ws.Range("J132").Value = 10624
$ws.Range("K132").Value = 22790.2851
$ws.Range("L132").Value = 31872
$ws.Range("M132").Value = -20260.2851
$ws.Range("N132").Value = -36932
$ws.Range("H136").Value = 4121.3057
$ws.Range("I136").Value = 2568.9524
$ws.Range("K136").Value = 7706.8572
$ws.Range("M136").Value = -5156.8572

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H24").Value = 29990
$ws.Range("I24").Value = 29990
$ws.Range("K24").Value = 29990
$ws.Range("M24").Value = -29760
$ws.Range("H33").Value = 26749.25
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 26749.25
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 26749.25
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -27249.25
$ws.Range("H36").Value = 26749.25
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 26749.25
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 26749.25
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -27249.25
$ws.Range("H43").Value = 60750
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 74333.336
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 74333.336
$ws.Range("M43").Value = -19851
$ws.Range("N43").Value = -74631.336
$ws.Range("H122").Value = 4201.3335
$ws.Range("I122").Value = 4280.421
$ws.Range("K122").Value = 12841.263
$ws.Range("M122").Value = -10391.263
$ws.Range("H130").Value = 55214.5
$ws.Range("J130").Value = 55214.5
$ws.Range("L130").Value = 55214.5
$ws.Range("N130").Value = -65254.5
$ws.Range("H135").Value = 71238
$ws.Range("J135").Value = 71238
$ws.Range("L135").Value = 71238
$ws.Range("N135").Value = -81378
